# Generate Report for Handoff
#
# A fresh handoff was generated for the "324a4567-38de-4d2d-975e-9d6a52a3674c.md"
# file (row 4 on every sheet). This stamps the new handoff/generation
# datetimes for that row on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for row 4
$wsOverview.Range("G4").Value = "2016-10-19 16:39:11"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for row 4
$wsZhCn.Range("H4").Value = "2016-10-19 16:38:59"

# de-de sheet: "Latest Handoff Datetime" column (H) for row 4
$wsDeDe.Range("H4").Value = "2016-10-19 16:39:11"
